$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 35: replace old "binary Fe-Ni taenite" entry with new
#     "Plagioclase / Sr (G25)" entry (Grocolas et al. 2025) ---
# Drop the old I35 cell entirely (no longer used for this row).
$ws.Range("I35").Clear()

$ws.Range("C35").Value = "Grocolas et al. (2025). Earth Planet. Sci. Lett."
$ws.Range("K35").Value = "PlSrG25"
$ws.Range("B35").Value = "Sr (G25)"
$ws.Range("A35").Value = "Plagioclase"
$ws.Range("D35").Value = -7.684
$ws.Range("F35").Value = -3.231
$ws.Range("G35").Value = -43640
$ws.Range("G35").NumberFormat = "0.00E+00"
$ws.Range("M35").Value = "lnfO2"
$ws.Range("N35").Value = "XAn"
$ws.Range("P35").Value = "P"
$ws.Range("Q35").Value = "P/T"
$ws.Range("R35").Value = "lnaSiO2"

# --- Row 36: brand-new row for "Plagioclase / Ba (G25)" entry (Grocolas et al. 2025) ---
$ws.Range("B36").Value = "Ba (G25)"
$ws.Range("K36").Value = "PlBaG25"
$ws.Range("A36").Value = "Plagioclase"
$ws.Range("C36").Value = "Grocolas et al. (2025). Earth Planet. Sci. Lett."
$ws.Range("D36").Value = -8.165
$ws.Range("F36").Value = -3.384
$ws.Range("G36").Value = -43490
$ws.Range("G36").NumberFormat = "0.00E+00"
$ws.Range("M36").Value = "lnfO2"
$ws.Range("N36").Value = "XAn"
$ws.Range("O36").Value = "1/T"
$ws.Range("P36").Value = "P"
$ws.Range("Q36").Value = "P/T"
$ws.Range("R36").Value = "lnaSiO2"

# --- Row 34: rename "Ba" -> "Ba (C)" and id "PlBa" -> "PlBaC" ---
$ws.Range("K34").Value = "PlBaC"
$ws.Range("B34").Value = "Ba (C)"

# --- Selection cursor moved by the author after editing ---
$ws.Range("A42").Select() | Out-Null
